$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pfizer block updates (rows 26-29) ---
$ws.Range("A26").Value = "5/31/2021"
$ws.Range("D26").Value = 'Delivered'

$ws.Range("A27").Value = "6/7/2021"
$ws.Range("D27").Value = 'Delivered'

$ws.Range("B28").Value = 165
$ws.Range("C28").Value = 193050
$ws.Range("D28").Value = 'Confirmed'

$ws.Range("B29").Value = 227
$ws.Range("C29").Value = 265590
$ws.Range("D29").Value = 'Confirmed'

# --- Moderna block updates (rows 47-50) ---
$ws.Range("D47").Value = 'Delivered'

$ws.Range("D48").Value = 'Delivered'

$ws.Range("A49").Value = "6/15/2021"
$ws.Range("D49").Value = 'Confirmed'

$ws.Range("A50").Value = "6/18/2021"
$ws.Range("D50").Value = 'Confirmed'

# --- AstraZeneca block updates (rows 73-80) ---
$ws.Range("B73").Value = 15840
$ws.Range("D73").Value = 'Delivered'

$ws.Range("A74").Value = "6/14/2021"
$ws.Range("B74").Value = 4800
$ws.Range("C74").Value = 48000
$ws.Range("D74").Value = 'Delivered'
$ws.Range("E74").Value = "6/14/2021"

$ws.Range("A75").Value = "6/22/2021"
$ws.Range("A76").Value = "6/29/2021"
$ws.Range("A77").Value = "7/6/2021"
$ws.Range("A78").Value = "7/13/2021"
$ws.Range("A79").Value = "7/20/2021"
$ws.Range("A80").Value = "7/27/2021"

# --- J&J block updates (rows 81-94, rows 81-82 converted from AstraZeneca) ---
$ws.Range("A81").Value = "4/14/2021"
$ws.Range("C81").Value = 9600
$ws.Range("D81").Value = 'Delivered'
$ws.Range("F81").Value = 'J&J'

$ws.Range("A82").Value = "4/29/2021"
$ws.Range("C82").Value = 6500
$ws.Range("D82").Value = 'Delivered'
$ws.Range("E82").Value = "4/30/2021"
$ws.Range("F82").Value = 'J&J'

$ws.Range("A83").Value = "5/14/2021"
$ws.Range("C83").Value = 12000
$ws.Range("E83").Value = "5/13/2021"

$ws.Range("A84").Value = "5/21/2021"
$ws.Range("C84").Value = 11600
$ws.Range("E84").Value = "5/19/2021"

$ws.Range("A85").Value = "5/28/2021"
$ws.Range("C85").Value = 26400
$ws.Range("E85").Value = "5/27/2021"

$ws.Range("A86").Value = "6/4/2021"
$ws.Range("C86").Value = 23760
$ws.Range("E86").Value = "5/28/2021"

$ws.Range("A87").Value = "6/11/2021"
$ws.Range("C87").Value = 12600
$ws.Range("E87").Value = "6/3/2021"

$ws.Range("A88").Value = "6/15/2021"
$ws.Range("C88").Value = 12000
$ws.Range("D88").Value = 'Confirmed'
$ws.Range("E88").Value = "6/7/2021"

$ws.Range("A89").Value = "6/22/2021"
$ws.Range("C89").Value = 23760

$ws.Range("A90").Value = "6/29/2021"
$ws.Range("C90").Value = 23760

$ws.Range("A91").Value = "7/6/2021"
$ws.Range("C91").Value = 23760

$ws.Range("A92").Value = "7/13/2021"
$ws.Range("C92").Value = 23760

$ws.Range("A93").Value = "7/20/2021"
$ws.Range("C93").Value = 23760

$ws.Range("A94").Value = "7/27/2021"
$ws.Range("C94").Value = 23760

# --- Row 95 converted from J&J Assumption row to a new Sputnik Delivered row ---
$ws.Range("A95").Value = "3/1/2021"
$ws.Range("C95").Value = 200000
$ws.Range("D95").Value = 'Delivered'
$ws.Range("F95").Value = 'Sputnik'

# --- Remove the old trailing row 96 (data no longer present) ---
$ws.Rows.Item(96).Delete()
